$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text even when it looks like a plain
# number (e.g. "212.02"), so it keeps matching the source inlineStr formatting
# (fixed decimals, no locale reformatting) instead of becoming a numeric cell.
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "26.665.90"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.596.81"
$ws.Range("E3").Value = "  -1.63%  "
Set-TextValue "D5" "212.02"
$ws.Range("E5").Value = "  -1.41%  "
Set-TextValue "D6" "0.515"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("E9").Value = "  -1.74%  "
Set-TextValue "D10" "19.72"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "1.820.51"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "1.590.44"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("E15").Value = "  -3.01%  "
Set-TextValue "D16" "65.11"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "26.639.63"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("E18").Value = "  -1.92%  "
Set-TextValue "D19" "210.37"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("E20").Value = "  +0.06%  "
Set-TextValue "D21" "6.72"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("E23").Value = "  -3.11%  "
Set-TextValue "D24" "8.89"
$ws.Range("E24").Value = "  -1.99%  "
Set-TextValue "D25" "146.74"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("E26").Value = "  +0.07%  "
Set-TextValue "D27" "7.17"
$ws.Range("E27").Value = "  -2.94%  "
Set-TextValue "D28" "0.116"
$ws.Range("E28").Value = "  -0.88%  "
Set-TextValue "D29" "15.34"
$ws.Range("E29").Value = "  -1.44%  "
Set-TextValue "D30" "0.0505"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("E31").Value = "  -1.42%  "
Set-TextValue "D32" "3.24"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("E33").Value = "  -10.00%  "
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("D35").Value = "1.292.83"
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -5.45%  "
$ws.Range("E38").Value = "  -3.45%  "
Set-TextValue "D39" "0.837"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "5.40"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "0.791"
$ws.Range("E42").Value = "  -1.13%  "
Set-TextValue "D43" "2.19"
$ws.Range("E43").Value = "  -1.53%  "
Set-TextValue "D44" "63.91"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").Value = "1.732.88"
$ws.Range("E45").Value = "  -1.66%  "
Set-TextValue "D46" "89.78"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  -2.39%  "
Set-TextValue "D48" "0.853"
$ws.Range("E48").Value = "  -3.24%  "
Set-TextValue "D49" "0.0985"
$ws.Range("E49").Value = "  -2.32%  "
Set-TextValue "D50" "0.0502"
$ws.Range("E50").Value = "  -2.36%  "
Set-TextValue "D51" "7.50"
$ws.Range("E51").Value = "  -2.25%  "
